# Hide (remove) the "Period" column from the "Details" sheet's UI.
# The "Period" column is column H on the "Details" sheet (between "Act #5"
# and "Debit"). Deleting the whole column shifts the following columns
# (Debit, Credit, Amount) left by one and removes the now-unused "Period"
# shared string entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

$ws.Columns.Item(8).Delete()
